# Generate Report for Handoff
# Rotate the localized-file UUID / content-hash that this handoff report
# references, and refresh the "generated at" timestamps that go with it.

$wb = $excel.ActiveWorkbook

$oldGuid = "23fe7231-6869-4b35-bbd6-8979f443ed0c"
$newGuid = "6c82ee76-023d-4634-b535-6d77c23f1aae"
$newHash = "1e98465918153e8e6845228c4b6ef77711dd9936"

# The hyperlinks all point at the same commit in the source repo; only the
# display text (file name) changes, the target URL itself is untouched.
$linkBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/03c05ca110dce3c2e7bdb16d33d0084a78c50ab6/e2e/" + $oldGuid + ".md"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview.Range("A2").Value = $newGuid + ".md"
$wsOverview.Range("B2").Value = "e2e\" + $newGuid + ".md"
$wsOverview.Range("G2").Value = "2016-08-28 00:55:14"

$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $linkBase, $null, $null, "e2e\" + $newGuid + ".md")

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn.Range("A2").Value = $newGuid + ".md"
$wsZhCn.Range("G2").Value = $newGuid + "." + $newHash + ".zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-08-28 00:55:10"

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $linkBase, $null, $null, $newGuid + ".md")

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe.Range("A2").Value = $newGuid + ".md"
$wsDeDe.Range("G2").Value = $newGuid + "." + $newHash + ".de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-08-28 00:55:14"

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $linkBase, $null, $null, $newGuid + ".md")
